$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before J (index 10). This shifts the old "Finance
# Related Issues" / "Location of Issue" columns (J, K) over to L, M, and
# relabels I from "Defects" to "Date of Last Update" while adding two fresh
# columns: a duplicate "TicketID" column and a new "Defects" column.
$ws.Range("J1:K1").EntireColumn.Insert()

$ws.Range("I1").Value = "Date of Last Update"

# Duplicate the TicketID column (header + value) via copy/paste so the new
# column keeps the same text-typed "248" entry instead of becoming numeric.
$ws.Range("A1:A2").Copy($ws.Range("J1:J2"))

$ws.Range("K1").Value = "Defects"
$ws.Range("K2").Value = "Incorrect / Error"

$wb.Save()
